# "setting warna ke ukuran" - rename the "warna" column header to "ukuran".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 held the "warna" header; change it to "ukuran".
$ws.Range("F1").Value = "ukuran"

# The author's last selection before saving moved to B8.
$ws.Range("B8").Select() | Out-Null
